$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "개인별내역" (Worksheets.Item(2)) - add 4 new attendance rows
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Fill in the previously-missing payment for 박상미 (row 7) so the
# outstanding-balance formula clears to 0 like the rest of the rows.
$ws2.Range("G7").Value = 10000

# --- copy formatting down for the 4 new rows (10-13) from row 9, which
#     already carries the styles used for the newer "카톡" entries ---
$ws2.Range("A9:H9").Copy()
$ws2.Range("A10:H13").PasteSpecial(-4122)
$ws2.Range("A9:H9").Copy()
$ws2.Range("A14:B14").PasteSpecial(-4122)
$ws2.Range("A9:H9").Copy()
$ws2.Range("D14:H14").PasteSpecial(-4122)
# column C of the new rows uses the header's style (s=36) rather than the
# "s=37" used by the rest of the data rows
$ws2.Range("C1").Copy()
$ws2.Range("C10:C13").PasteSpecial(-4122)

# --- row 10: 박진규 / 4조 / 지각 ---
$ws2.Range("A10").Value = "박진규"
$ws2.Range("B10").Formula = "=VLOOKUP(A10,주소록!`$A`$2:`$D`$28,4,FALSE)"
$ws2.Range("C10").Value = "지각"
$ws2.Range("D10").Value = 42404
$ws2.Range("G10").Value = 5000

# --- row 11: 윤승업 / 3조 / 지각 ---
$ws2.Range("A11").Value = "윤승업"
$ws2.Range("B11").Formula = "=VLOOKUP(A11,주소록!`$A`$2:`$D`$28,4,FALSE)"
$ws2.Range("C11").Value = "지각"
$ws2.Range("D11").Value = 42415
$ws2.Range("G11").Value = 10000

# --- row 12: 윤성민 / 5조 / 지각 ---
$ws2.Range("A12").Value = "윤성민"
$ws2.Range("B12").Formula = "=VLOOKUP(A12,주소록!`$A`$2:`$D`$28,4,FALSE)"
$ws2.Range("C12").Value = "지각"
$ws2.Range("D12").Value = 42415
$ws2.Range("G12").Value = 10000

# --- row 13: 유진혁 / 4조 / 지각 ---
$ws2.Range("A13").Value = "유진혁"
$ws2.Range("B13").Formula = "=VLOOKUP(A13,주소록!`$A`$2:`$D`$28,4,FALSE)"
$ws2.Range("C13").Value = "지각"
$ws2.Range("D13").Value = 42416
$ws2.Range("G13").Value = 5000

# --- shared formulas for the whole column so Excel regroups them the
#     same way the authoring session did ---
$ws2.Range("E2:E13").Formula = "=IF(D2<>0,VLOOKUP(WEEKDAY(D2,2),weekday,2,FALSE),`"`")"
$ws2.Range("F2:F13").Formula = "=IF(C2<>0,IF(E2=`"월`",VLOOKUP(C2,cost,2,FALSE)*2,VLOOKUP(C2,cost,2,FALSE)),`"`")"
$ws2.Range("H2:H6").Formula = "=F2-G2"
$ws2.Range("H7:H13").Formula = "=F7-G7"

# row 14 stays blank (just formatted) - matches the trailing blank table row
$ws2.Range("A14:H14").ClearContents()

# extend the "data" table to include the new rows
$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:H14"))

# data-validation list on column C now spans the new rows too
$ws2.Range("C2:C13").Validation.Delete()
$ws2.Range("C2:C13").Validation.Add(3, 1, 1, "=type")

# printed page size, matching the rest of the workbook
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws2.Range("A14").Select()

# ---------------------------------------------------------------------
# Sheet "지출내용" (Worksheets.Item(3)) - add the MT deposit refund row
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B2").Copy()
$ws3.Range("A4").PasteSpecial(-4122)
$ws3.Range("B3").Copy()
$ws3.Range("B4:C4").PasteSpecial(-4122)

$ws3.Range("A4").Value = "MT비 보충"
$ws3.Range("B4").Value = 42412
$ws3.Range("C4").Value = 59000

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:C4"))

$ws3.Range("D5").Select()
$ws3.Activate()

# ---------------------------------------------------------------------
# Workbook window bookkeeping
# ---------------------------------------------------------------------
$wb.Save()
